$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated vm_pu results for the 380 kV case (Case_5_106).
# Column order per row: B,C,D,E,F,I,J,K,L,M (G,H,N unchanged).
$newValues = @{
    "B2" = 1.02; "C2" = 1.071720654839941; "D2" = 1.07233951787018; "E2" = 1.07553661343965; "F2" = 1.08489802213939; "I2" = 1.056687505477585; "J2" = 1.076643822910361; "K2" = 1.075034294231136; "L2" = 1.078222923148158; "M2" = 1.087559850481639
    "B3" = 1.02; "C3" = 1.072846994044983; "D3" = 1.07321474526172; "E3" = 1.076516678603261; "F3" = 1.08590222190778; "I3" = 1.057023721535139; "J3" = 1.077427355294652; "K3" = 1.075726102761785; "L3" = 1.079019922956136; "M3" = 1.088382694705964
    "B4" = 1.02; "C4" = 1.073576123977232; "D4" = 1.073781312809024; "E4" = 1.077151391187133; "F4" = 1.086552567022073; "I4" = 1.057240263459715; "J4" = 1.077934059826055; "K4" = 1.076173343002859; "L4" = 1.079535557031643; "M4" = 1.088915070758559
    "B5" = 1.02; "C5" = 1.073882725371217; "D5" = 1.07401955402564; "E5" = 1.077418354377002; "F5" = 1.086826105889157; "I5" = 1.057331055177647; "J5" = 1.078147007994584; "K5" = 1.076361265291357; "L5" = 1.079752310608274; "M5" = 1.089138867148556
    "B6" = 1.02; "C6" = 1.073934209508419; "D6" = 1.074059559056642; "E6" = 1.07746318628332; "F6" = 1.086872042100385; "I6" = 1.057346285281576; "J6" = 1.078182758794213; "K6" = 1.076392812544494; "L6" = 1.0797887033656; "M6" = 1.08917644268711
    "B7" = 1.02; "C7" = 1.073580220504426; "D7" = 1.073784495979531; "E7" = 1.077154957851983; "F7" = 1.086556221536698; "I7" = 1.057241477576161; "J7" = 1.077936905526832; "K7" = 1.076175854413162; "L7" = 1.079538453378516; "M7" = 1.088918061194273
    "B8" = 1.02; "C8" = 1.072101241727774; "D8" = 1.07263525571553; "E8" = 1.075867717970112; "F8" = 1.085237279778728; "I8" = 1.056801340791072; "J8" = 1.076908681962087; "K8" = 1.075268178088698; "L8" = 1.078492288958121; "M8" = 1.087837946592531
    "B9" = 1.02; "C9" = 1.069497481231969; "D9" = 1.070611983483146; "E9" = 1.073603628622904; "F9" = 1.082917450115906; "I9" = 1.05601801666587; "J9" = 1.075094582635956; "K9" = 1.073665637970375; "L9" = 1.076648228150465; "M9" = 1.085934207889327
    "B10" = 1.02; "C10" = 1.067763231388565; "D10" = 1.069264396814444; "E10" = 1.072097076903654; "F10" = 1.081373822127672; "I10" = 1.05549059966333; "J10" = 1.07388368753645; "K10" = 1.072595209072355; "L10" = 1.075418475714622; "M10" = 1.084664771230003
    "B11" = 1.02; "C11" = 1.067012654893675; "D11" = 1.068681181181303; "E11" = 1.071445401361244; "F11" = 1.080706112393262; "I11" = 1.055260989920944; "J11" = 1.073359003086668; "K11" = 1.072131213267674; "L11" = 1.074885892044215; "M11" = 1.084115028323797
    "B12" = 1.02; "C12" = 1.066733911733052; "D12" = 1.068464594071114; "E12" = 1.071203441148125; "F12" = 1.08045819940926; "I12" = 1.055175517178846; "J12" = 1.073164058007953; "K12" = 1.071958790607647; "L12" = 1.074688052937296; "M12" = 1.083910819202641
    "B13" = 1.02; "C13" = 1.066793700634688; "D13" = 1.068511050693327; "E13" = 1.071255337856292; "F13" = 1.080511372868922; "I13" = 1.055193859763984; "J13" = 1.073205876828633; "K13" = 1.071995779198474; "L13" = 1.074730490720734; "M13" = 1.083954623208137
    "B14" = 1.02; "C14" = 1.066989612790753; "D14" = 1.068663277100383; "E14" = 1.071425398790731; "F14" = 1.080685617702647; "I14" = 1.055253928500516; "J14" = 1.073342889965405; "K14" = 1.072116962264881; "L14" = 1.074869538883131; "M14" = 1.084098148539844
    "B15" = 1.02; "C15" = 1.067110327925951; "D15" = 1.068757074805723; "E15" = 1.071530192349144; "F15" = 1.080792989503351; "I15" = 1.055290914257538; "J15" = 1.073427301137261; "K15" = 1.072191617354817; "L15" = 1.074955209217598; "M15" = 1.084186577896237
    "B16" = 1.02; "C16" = 1.067813052329692; "D16" = 1.06930310922789; "E16" = 1.072140340627406; "F16" = 1.081418150413124; "I16" = 1.055505812095309; "J16" = 1.073918501554546; "K16" = 1.072625992567841; "L16" = 1.075453819594524; "M16" = 1.084701254427145
    "B17" = 1.02; "C17" = 1.068253950192573; "D17" = 1.069645702287056; "E17" = 1.072523250419282; "F17" = 1.081810482625368; "I17" = 1.055640281294107; "J17" = 1.074226522519276; "K17" = 1.072898332744541; "L17" = 1.075766559960377; "M17" = 1.08502407919953
    "B18" = 1.02; "C18" = 1.068511153718934; "D18" = 1.069845559780494; "E18" = 1.072746659976168; "F18" = 1.082039390056073; "I18" = 1.055718595769302; "J18" = 1.074406151309375; "K18" = 1.073057136584207; "L18" = 1.075948967167629; "M18" = 1.085212370811921
    "B19" = 1.02; "C19" = 1.068598859470714; "D19" = 1.069913710888865; "E19" = 1.072822847837553; "F19" = 1.082117452877857; "I19" = 1.055745278773268; "J19" = 1.07446739421126; "K19" = 1.07311127651625; "L19" = 1.076011161788709; "M19" = 1.085276572259018
    "B20" = 1.02; "C20" = 1.068206642380543; "D20" = 1.069608942317018; "E20" = 1.072482161122809; "F20" = 1.081768382191987; "I20" = 1.055625866343882; "J20" = 1.074193478366383; "K20" = 1.072869118141283; "L20" = 1.075733006804408; "M20" = 1.084989443839706
    "B21" = 1.02; "C21" = 1.066931920046811; "D21" = 1.06861844896881; "E21" = 1.071375317291386; "F21" = 1.080634304065064; "I21" = 1.055236244874066; "J21" = 1.073302544499777; "K21" = 1.072081278912933; "L21" = 1.074828593045498; "M21" = 1.084055884185259
    "B22" = 1.02; "C22" = 1.066130765846992; "D22" = 1.067995947839466; "E22" = 1.070679985695509; "F22" = 1.07992186679921; "I22" = 1.054990201396907; "J22" = 1.072742067100719; "K22" = 1.071585505276173; "L22" = 1.074259872329378; "M22" = 1.083468859656525
    "B23" = 1.02; "C23" = 1.066555443213757; "D23" = 1.068325922582863; "E23" = 1.07104853859222; "F23" = 1.080299486057189; "I23" = 1.055120735393103; "J23" = 1.073039216338018; "K23" = 1.071848364777689; "L23" = 1.074561369474736; "M23" = 1.083780057960426
    "B24" = 1.02; "C24" = 1.068228018635275; "D24" = 1.069625552476248; "E24" = 1.07250072740706; "F24" = 1.081787405359505; "I24" = 1.055632380207449; "J24" = 1.074208409703586; "K24" = 1.07288231910989; "L24" = 1.075748168058896; "M24" = 1.085005094088302
    "B25" = 1.02; "C25" = 1.070170334259381; "D25" = 1.071134827494607; "E25" = 1.074188450890803; "F25" = 1.083516667915173; "I25" = 1.056221441994904; "J25" = 1.075563835077408; "K25" = 1.074080298709716; "L25" = 1.07712502990455; "M25" = 1.086426420179064
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value2 = $newValues[$addr]
}
